$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# Row 41 - new bug entry
$ws.Range("A41").Value = "SB"
$ws.Range("B41").Value = 43664
$ws.Range("B2").Copy()
$ws.Range("B41").PasteSpecial($xlPasteFormats)
$ws.Range("C41").Value = "confusion caused by the meta_table object in the cctu library, versus the local users object which might be called meta_table"

# Row 42 - new bug entry
$ws.Range("A42").Value = "SB"
$ws.Range("B42").Value = 43664
$ws.Range("B2").Copy()
$ws.Range("B42").PasteSpecial($xlPasteFormats)
$ws.Range("C42").Value = 'make it easier to get at the cctu:::cctu_env$code_tree object, or at least document it better'

$excel.CutCopyMode = $false

# Update the view to match the new selection / scroll position
$ws.Range("B42").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 2
